# Update the cryptos list with refreshed prices and 1h volume percentages
# (GitHub Actions scheduled data refresh).
#
# Cell values are prefixed with a leading apostrophe to force text entry
# (mirrors how Excel's UI treats a leading "'" ) so purely-numeric-looking
# strings such as "570.27" or "6.25" stay text cells (matching the source
# workbook's inlineStr cells) instead of being auto-coerced to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.150.09"
$ws.Range("E2").Value = "'  +0.33%  "
$ws.Range("D3").Value = "'2.573.83"
$ws.Range("E3").Value = "'  -1.22%  "
$ws.Range("E4").Value = "'  -0.01%  "
$ws.Range("D5").Value = "'570.27"
$ws.Range("E5").Value = "'  +2.55%  "
$ws.Range("D6").Value = "'143.32"
$ws.Range("E6").Value = "'  -0.42%  "
$ws.Range("E7").Value = "'  +0.13%  "
$ws.Range("D8").Value = "'0.596"
$ws.Range("E8").Value = "'  +0.27%  "
$ws.Range("D9").Value = "'2.580.22"
$ws.Range("E9").Value = "'  -1.45%  "
$ws.Range("E10").Value = "'  -1.48%  "
$ws.Range("D11").Value = "'0.104"
$ws.Range("E11").Value = "'  +2.97%  "
$ws.Range("E12").Value = "'  +11.70%  "
$ws.Range("D13").Value = "'0.346"
$ws.Range("E13").Value = "'  +3.11%  "
$ws.Range("D14").Value = "'3.020.87"
$ws.Range("E14").Value = "'  -1.45%  "
$ws.Range("D15").Value = "'59.166.97"
$ws.Range("E15").Value = "'  +0.46%  "
$ws.Range("D16").Value = "'22.50"
$ws.Range("E16").Value = "'  +7.72%  "
$ws.Range("E17").Value = "'  +3.65%  "
$ws.Range("D18").Value = "'2.578.05"
$ws.Range("E18").Value = "'  -1.40%  "
$ws.Range("D19").Value = "'4.54"
$ws.Range("E19").Value = "'  +1.79%  "
$ws.Range("D20").Value = "'337.89"
$ws.Range("E20").Value = "'  +0.08%  "
$ws.Range("D21").Value = "'10.24"
$ws.Range("E21").Value = "'  +1.31%  "
$ws.Range("D22").Value = "'6.25"
$ws.Range("E22").Value = "'  +1.17%  "
$ws.Range("D24").Value = "'64.46"
$ws.Range("E24").Value = "'  -3.16%  "
$ws.Range("D25").Value = "'0.458"
$ws.Range("E25").Value = "'  +6.76%  "
$ws.Range("D26").Value = "'0.996"
$ws.Range("E26").Value = "'  -0.16%  "
$ws.Range("E27").Value = "'  +0.15%  "
$ws.Range("D28").Value = "'7.26"
$ws.Range("E28").Value = "'  +0.97%  "
$ws.Range("D29").Value = "'0.0₃0781"
$ws.Range("E29").Value = "'  +2.93%  "
$ws.Range("E30").Value = "'  +0.03%  "
$ws.Range("E31").Value = "'  +0.10%  "
$ws.Range("E32").Value = "'  +1.69%  "
$ws.Range("D33").Value = "'158.52"
$ws.Range("E33").Value = "'  +2.86%  "
$ws.Range("D34").Value = "'19.02"
$ws.Range("E34").Value = "'  +0.01%  "
$ws.Range("D35").Value = "'4.04"
$ws.Range("E35").Value = "'  +2.26%  "
$ws.Range("E36").Value = "'  +1.84%  "
$ws.Range("D37").Value = "'0.872"
$ws.Range("E37").Value = "'  -3.52%  "
$ws.Range("D38").Value = "'0.872"
$ws.Range("E38").Value = "'  -1.47%  "
$ws.Range("D39").Value = "'37.20"
$ws.Range("E40").Value = "'  +2.64%  "
$ws.Range("E41").Value = "'  +2.08%  "
$ws.Range("D42").Value = "'293.23"
$ws.Range("E42").Value = "'  +3.86%  "
$ws.Range("E43").Value = "'  +0.24%  "
$ws.Range("D44").Value = "'0.0977"
$ws.Range("E44").Value = "'  +2.50%  "
$ws.Range("D45").Value = "'0.592"
$ws.Range("E45").Value = "'  -1.35%  "
# Rows 46-49 were reshuffled: Aave moved up to row 46, pushing Hedera,
# EnergySwap and WhiteBITCoin down by one row each (with refreshed values).
$ws.Range("B46").Value = "'Aave"
$ws.Range("C46").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'127.37"
$ws.Range("E46").Value = "'  +6.96%  "
$ws.Range("B47").Value = "'Hedera"
$ws.Range("C47").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").Value = "'0.0536"
$ws.Range("E47").Value = "'  -0.14%  "
$ws.Range("B48").Value = "'EnergySwap"
$ws.Range("C48").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'19.23"
$ws.Range("E48").Value = "'  +2.12%  "
$ws.Range("B49").Value = "'WhiteBITCoin"
$ws.Range("C49").Value = "'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").Value = "'10.64"
$ws.Range("E49").Value = "'  +0.22%  "
$ws.Range("E50").Value = "'  +1.95%  "
$ws.Range("D51").Value = "'1.948.41"
